$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.035.01"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "1.849.22"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "237.77"
$ws.Range("E5").Value = "  +3.37%  "
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "42.42"
$ws.Range("E8").Value = "  +6.46%  "
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D11").Value = "0.0990"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "2.114.44"
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "11.37"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.850.75"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D17").Value = "34.998.22"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").Value = "69.95"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "0.0₃0791"
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("D20").Value = "240.35"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "12.15"
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("D22").Value = "4.74"
$ws.Range("E22").Value = "  +1.91%  "
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").Value = "169.55"
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "1.85"
$ws.Range("E26").Value = "  +22.74%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "7.97"
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("D28").Value = "17.58"
$ws.Range("E28").Value = "  +1.79%  "
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").Value = "0.0554"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("D34").Value = "1.67"
$ws.Range("E34").Value = "  +23.46%  "
$ws.Range("E35").Value = "  +9.40%  "
$ws.Range("E36").Value = "  +4.95%  "
$ws.Range("D37").Value = "0.783"
$ws.Range("E37").Value = "  +13.83%  "
$ws.Range("B38").Value = "Gas"
$ws.Range("C38").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D38").Value = "14.28"
$ws.Range("E38").Value = "  +67.52%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "1.07"
$ws.Range("E39").Value = "  +10.39%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.0203"
$ws.Range("E40").Value = "  +5.66%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "90.17"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.342.05"
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("D43").Value = "14.90"
$ws.Range("E43").Value = "  +5.01%  "
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("E46").Value = "  +6.00%  "
$ws.Range("D47").Value = "2.73"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").Value = "6.48"
$ws.Range("E48").Value = "  +4.96%  "
$ws.Range("D49").Value = "2.028.06"
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("D50").Value = "0.0674"
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("E51").Value = "  +0.11%  "
